$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.443.41"
$ws.Range("E2").Value = "  +5.81%  "

$ws.Range("D3").Value = "2.047.49"
$ws.Range("E3").Value = "  +3.15%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.647"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "66.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +17.13%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").Value = "0.376"
$ws.Range("E10").Value = "  +5.40%  "

$ws.Range("E11").Value = "  +4.24%  "

$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.86%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.909"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("D15").Value = "2.354.43"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.12%  "

$ws.Range("D17").Value = "20.68"
$ws.Range("E17").Value = "  +20.67%  "

$ws.Range("D18").Value = "2.078.99"
$ws.Range("E18").Value = "  +4.81%  "

$ws.Range("D19").Value = "37.404.72"
$ws.Range("E19").Value = "  +5.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.17%  "

$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  +5.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.89%  "

$ws.Range("D23").Value = "237.56"
$ws.Range("E23").Value = "  +2.71%  "

$ws.Range("D24").Value = "2.75"
$ws.Range("E24").Value = "  +22.25%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.98%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.65%  "

$ws.Range("E30").Value = "  +10.38%  "

$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +3.03%  "

$ws.Range("E32").Value = "  +7.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +23.31%  "

$ws.Range("D34").Value = "4.71"
$ws.Range("E34").Value = "  +10.96%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0614"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.33%  "

$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +9.14%  "

$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  +24.78%  "

$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("E40").Value = "  +17.72%  "

$ws.Range("E41").Value = "  +4.77%  "

$ws.Range("E42").Value = "  +5.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0219"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.76%  "

$ws.Range("E44").Value = "  +6.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +20.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.19%  "

$ws.Range("D49").Value = "1.429.89"
$ws.Range("E49").Value = "  +4.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.28%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.38%  "
